{"js": "// Supplementary Figure 3 caption: \"Scree plot on sd rank correlation matrix.\"\n// becomes \"Scree plot on mean rank correlation matrix.\" (replace the word\n// \"sd\" with \"mean\", keeping the existing bold Times New Roman formatting).\nconst body = context.document.body;\n\nconst results = body.search(\"sd\", { matchCase: true, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the word \"sd\" to replace.');\n}\n\n// The caption text is unique in the document (\"sd\" only occurs as a whole\n// word inside \"Scree plot on sd rank correlation matrix.\"), so this loop\n// touches exactly that one occurrence; insertText(\"Replace\") swaps just\n// the matched word's text while leaving its run's formatting untouched.\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"mean\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Supplementary Figure 3 caption: \"Scree plot on sd rank correlation matrix.\"\n# becomes \"Scree plot on mean rank correlation matrix.\" (replace the word\n# \"sd\" with \"mean\", keeping the existing bold Times New Roman formatting).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"sd\"\n$rng.Find.MatchWholeWord = $true\n$rng.Find.MatchCase = $true\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"mean\"\n\n# wdReplaceOne (1): only the single \"sd\" caption word qualifies as a whole-\n# word, case-sensitive match, so this touches just that occurrence while\n# keeping the run's existing (bold, Times New Roman) character formatting.\n$found = $rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\nif (-not $found) {\n    throw 'Could not find the word \"sd\" to replace.'\n}\n"}
